$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the missing Friday hours for the week of row 8
$ws.Range("F8").Value = 7.25

# Move the active selection to reflect where the user left off editing
$ws.Range("H8").Select()

$wb.Save()
